$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update language code for the row from English to French
$ws.Range("A2").Value = "fra"

# Update the translated name/description to French for both the name and descr columns
$ws.Range("C2").Value = "Machine virtuelle rÃ©sidente"
$ws.Range("H2").Value = "Machine virtuelle rÃ©sidente"

# Drop the custom boolean number format previously applied to the is_active cell
$ws.Range("I2").ClearFormats()
